# Applies the "updated xls, added RNN file" commit to the workbook.
# Target sheet: "Advanced Learning"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Advanced Learning")

# ---------------------------------------------------------------------
# 1. New "Deep Networks" section header label (A45) and renamed model
#    labels in A46:A51 (new deep-learning model runs replacing the old
#    placeholder labels that were copy/pasted from the other sections).
#    NOTE: assigned in this particular order so the shared-string table
#    is built up in the same order as the authoritative workbook.
# ---------------------------------------------------------------------
$ws.Range("A45").Value = "Deep Networks"
$ws.Range("A51").Value = "250-100-10-7"
$ws.Range("A49").Value = "250-100-10-7_32-16"
$ws.Range("A46").Value = "250-100-10-7_1-64"
$ws.Range("A47").Value = "250-100-10-7_1-64-bin"
$ws.Range("A48").Value = "250-100-10-7_1-64-mult"
$ws.Range("A50").Value = "250-100-10-7_RNN"

# ---------------------------------------------------------------------
# 2. New score rows for the Deep Networks runs (Accuracy/Precision/
#    Recall/F1 scores per dataset variant).
# ---------------------------------------------------------------------
$ws.Range("B46").Formula = "=AVERAGE(0.83687909, 0.81979017, 0.96844162, 0.99526694, 0.98501774, 0.97060373, 0.98451186)"
$ws.Range("C46").Value = 0.78
$ws.Range("D46").Value = 0.77
$ws.Range("E46").Value = 0.76

$ws.Range("B47").Formula = "=AVERAGE(0.72299451, 0.69590587,0.95264858, 0.99526694,0.98354169, 0.97014636, 0.97528135)"
$ws.Range("C47").Value = 0.62
$ws.Range("D47").Value = 0.63
$ws.Range("E47").Value = 0.61

$ws.Range("B48").Formula = "=AVERAGE(0.77954873, 0.76406752, 0.95583629, 0.99526694,0.98354169,0.97014636, 0.97384688)"
$ws.Range("C48").Value = 0.69
$ws.Range("D48").Value = 0.7
$ws.Range("E48").Value = 0.69

$ws.Range("B49").Formula = "=AVERAGE(0.92409774,0.91172109,0.98151125,0.99745676,0.99081107,0.98283485, 0.9908665)"
$ws.Range("C49").Value = 0.9
$ws.Range("D49").Value = 0.88
$ws.Range("E49").Value = 0.89

# ---------------------------------------------------------------------
# 3. Newly measured scores for the "Numeric Only min-max 1" (rows
#    30-32) dataset variant that previously had no results.
# ---------------------------------------------------------------------
$ws.Range("B30").Value = 0.65996091584432803
$ws.Range("C30").Value = 0.63
$ws.Range("D30").Value = 0.66
$ws.Range("E30").Value = 0.63

$ws.Range("B31").Value = 0.48720756181394798
$ws.Range("C31").Value = 0.24
$ws.Range("D31").Value = 0.49
$ws.Range("E31").Value = 0.32

$ws.Range("B32").Value = 0.48079055327641601
$ws.Range("C32").Value = 0.37
$ws.Range("D32").Value = 0.48
$ws.Range("E32").Value = 0.33

# ---------------------------------------------------------------------
# 4. Updated row-count (Support) figure: 120253 -> 144304, now shown in
#    bold, across every dataset-variant block that references it.
# ---------------------------------------------------------------------
$supportRows = 22,23,24,25,26,27,30,31,32,33,34,35,38,39,40,41,42,43
foreach ($r in $supportRows) {
    $cell = $ws.Range("F$r")
    $cell.Value = 144304
    $cell.Font.Bold = $true
}

# ---------------------------------------------------------------------
# 5. Restore the view state (scroll position / active selection) that
#    was in effect when the workbook was last saved.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A50").Select()

# ---------------------------------------------------------------------
# 6. Recalculate so every dependent formula (same sheet + "ALL
#    Learning") picks up the new cached values.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()
